# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The workbook holds a small "Estado de Cuenta" (account statement) table in
# Hoja1!B15:J22. Rows 16-22 list one worker-arrears record per row (doc type,
# doc number, name, arrears period, arrears value, base salary). The database
# backing the sheet was refreshed: the old arrears records are dropped and the
# new ones are written in, so every data row (16-22) ends up re-populated -
# same columns/styles, new values, and the rows now come out sorted by
# worker document number.
#
# Values are written row-by-row, in the exact order they appear in the
# refreshed table, so that C (doc #) -> D (name) -> E (period) -> F (valor
# mora) -> G (salario basico) are all touched per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16: ENOTH ENRIQUE GARCIA YEPEZ
$ws.Cells.Item(16, 3).Value = "73214033"
$ws.Cells.Item(16, 4).Value = "ENOTH ENRIQUE GARCIA YEPEZ"
$ws.Cells.Item(16, 5).Value = "2006"
$ws.Cells.Item(16, 6).Value = 35112
$ws.Cells.Item(16, 7).Value = 877803

# Row 17: JOHNNY FERNANDO REINA BOLIVAR
$ws.Cells.Item(17, 3).Value = "1047421288"
$ws.Cells.Item(17, 4).Value = "JOHNNY FERNANDO REINA BOLIVAR"
$ws.Cells.Item(17, 5).Value = "1903"
$ws.Cells.Item(17, 6).Value = 4417
$ws.Cells.Item(17, 7).Value = 828116

# Row 18: HAROLD JEYSON HERRERA SAYAVEDRA
$ws.Cells.Item(18, 3).Value = "73205202"
$ws.Cells.Item(18, 4).Value = "HAROLD JEYSON HERRERA SAYAVEDRA"
$ws.Cells.Item(18, 5).Value = "1905"
$ws.Cells.Item(18, 6).Value = 1893
$ws.Cells.Item(18, 7).Value = 1420000

# Row 19: CARLOS DE JESUS MARTINEZ VILORIA
$ws.Cells.Item(19, 3).Value = "1143363534"
$ws.Cells.Item(19, 4).Value = "CARLOS DE JESUS MARTINEZ VILORIA"
$ws.Cells.Item(19, 5).Value = "2412"
$ws.Cells.Item(19, 6).Value = 24800
$ws.Cells.Item(19, 7).Value = 1550000

# Row 20: EDER JULIAN ARCHBOLD SALCEDO
$ws.Cells.Item(20, 3).Value = "1128046927"
$ws.Cells.Item(20, 4).Value = "EDER JULIAN ARCHBOLD SALCEDO"
$ws.Cells.Item(20, 5).Value = "1809"
$ws.Cells.Item(20, 6).Value = 1200
$ws.Cells.Item(20, 7).Value = 900000

# Row 21: ALEMIS VILLARREAL ANGULO
$ws.Cells.Item(21, 3).Value = "1143347543"
$ws.Cells.Item(21, 4).Value = "ALEMIS VILLARREAL ANGULO"
$ws.Cells.Item(21, 5).Value = "1809"
$ws.Cells.Item(21, 6).Value = 1053
$ws.Cells.Item(21, 7).Value = 790000

# Row 22: VERONICA RAMOS LARA
$ws.Cells.Item(22, 3).Value = "1143401657"
$ws.Cells.Item(22, 4).Value = "VERONICA RAMOS LARA"
$ws.Cells.Item(22, 5).Value = "1905"
$ws.Cells.Item(22, 6).Value = 1104
$ws.Cells.Item(22, 7).Value = 877803
